# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Leading "'" on numeric-looking Price values keeps them stored as text
# (matching the source data), since Excel would otherwise coerce a bare
# "247.49"-style string into a numeric cell on assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.177.73"
$ws.Range("E2").Value = "  +1.37%  "

$ws.Range("D3").Value = "2.023.76"
$ws.Range("E3").Value = "  +3.26%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'247.49"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("E6").Value = "  +1.93%  "

$ws.Range("D7").Value = "'60.16"
$ws.Range("E7").Value = "  -2.66%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.393"
$ws.Range("E9").Value = "  +4.46%  "

$ws.Range("D10").Value = "'0.0812"
$ws.Range("E10").Value = "  +2.36%  "

$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("E12").Value = "  +6.65%  "

$ws.Range("D13").Value = "'22.55"
$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").Value = "'0.858"
$ws.Range("E14").Value = "  +3.05%  "

$ws.Range("D15").Value = "2.317.71"
$ws.Range("E15").Value = "  +3.26%  "

$ws.Range("D16").Value = "'5.51"
$ws.Range("E16").Value = "  +4.01%  "

$ws.Range("D17").Value = "2.021.44"
$ws.Range("E17").Value = "  +3.36%  "

$ws.Range("D18").Value = "37.137.62"
$ws.Range("E18").Value = "  +1.61%  "

$ws.Range("D19").Value = "'70.69"
$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("E21").Value = "  +3.42%  "

$ws.Range("D22").Value = "'231.01"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'2.52"
$ws.Range("E24").Value = "  +2.58%  "

$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").Value = "'9.50"
$ws.Range("E26").Value = "  +3.35%  "

$ws.Range("D27").Value = "'163.58"
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("D28").Value = "'0.138"
$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("D29").Value = "'19.83"
$ws.Range("E29").Value = "  +2.10%  "

$ws.Range("E30").Value = "  +6.75%  "

$ws.Range("E31").Value = "  +2.19%  "

$ws.Range("E32").Value = "  +1.33%  "

$ws.Range("D33").Value = "'0.0664"
$ws.Range("E33").Value = "  +7.86%  "

$ws.Range("D34").Value = "'4.54"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("E35").Value = "  +8.25%  "

$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("E38").Value = "  +1.57%  "

$ws.Range("D39").Value = "'5.45"
$ws.Range("E39").Value = "  -1.28%  "

$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("D44").Value = "'16.76"
$ws.Range("E44").Value = "  +4.53%  "

$ws.Range("D45").Value = "'92.22"
$ws.Range("E45").Value = "  +3.97%  "

$ws.Range("D46").Value = "1.392.86"
$ws.Range("E46").Value = "  +1.69%  "

$ws.Range("E47").Value = "  +2.99%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.48"
$ws.Range("E48").Value = "  +4.74%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.19"
$ws.Range("E49").Value = "  +18.41%  "

$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").Value = "'46.93"
$ws.Range("E51").Value = "  +3.56%  "
